# Applies the cryptos-list refresh described by the commit:
# "Updated cryptos list ... with GitHub Actions" - new Price (D) and
# Volume(1h) (E) figures for every coin row, plus a few re-ranked coins
# (B/C pairs swapped between adjacent rows 41/42 and 45/46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.752.53"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.877.80"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'332.88"
$ws.Range("E5").Value = "  +3.67%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4725"
$ws.Range("E7").Value = "  +6.39%  "
$ws.Range("D8").Value = "'0.3968"
$ws.Range("E8").Value = "  +4.18%  "
$ws.Range("D9").Value = "'47.94"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "'0.08038"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("D11").Value = "'1.026"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "'21.91"
$ws.Range("E12").Value = "  +2.95%  "
$ws.Range("D13").Value = "1.882.97"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("D15").Value = "'7.176"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "'0.00001052"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "'87.26"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").Value = "'0.06628"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "'17.35"
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "27.791.97"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "'5.504"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").Value = "'11.06"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("D25").Value = "'2.299"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "2.109.57"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").Value = "'156.80"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").Value = "'20.23"
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("D29").Value = "'2.105"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "'5.606"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").Value = "'122.59"
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("E32").Value = "  +5.49%  "
$ws.Range("D33").Value = "'0.09565"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "'1.458"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").Value = "'5.309"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "'0.02268"
$ws.Range("E37").Value = "  +2.77%  "
$ws.Range("D38").Value = "'0.06113"
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("D39").Value = "'1.232"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").Value = "'8.185"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6011"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "'1.002"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'0.1914"
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("D44").Value = "'10.30"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.263"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5711"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").Value = "'12.26"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "'3.413"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "'1.941"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").Value = "'0.06808"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'113.35"
$ws.Range("E51").Value = "  +5.22%  "
